$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix the "TENCHUONGTRINH" (program name) text in column G: drop the stray
# leading space so the shared string reads "CONG NGHE THONG TIN KHOA 2020"
# instead of " CONG NGHE THONG TIN KHOA 2020".
$ws.Range("G3:G68").Value = "CÔNG NGHỆ THÔNG TIN KHÓA 2020"

# Move the viewport/selection to where the user was last working (around the
# newly added "lop" rows), mirroring the cursor move captured in the saved
# sheet view.
$excel.ActiveWindow.ScrollRow = 21
$excel.ActiveWindow.ScrollColumn = 2
$ws.Range("G45").Select()
